$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The target cells hold numeric-looking values that are stored as TEXT
# (shared strings) in the workbook, not as real numbers. Writing a plain
# numeric string via .Value would make Excel coerce the cell to a Number
# type, which doesn't match the source data. Prefixing with a leading
# apostrophe forces Excel to keep/store the value as text (quoted text),
# matching the original "number-as-text" representation.
#
# That text coercion also causes Excel to tag the cell with a new
# "quote prefix" cell style, even though the visual style is identical to
# before. We preserve the original Style object and re-apply it after the
# write so the cell's style reference is left unchanged.

function Set-TextValue($cellRef, $newValue) {
    $cell = $ws.Range($cellRef)
    $originalStyle = $cell.Style
    $cell.Value = "'" + $newValue
    $cell.Style = $originalStyle
}

# Row 11 - Enterprises density (per 1000 people)
Set-TextValue "B11" "10.23"
Set-TextValue "C11" "4.17"

# Row 12 - Employment (% of total)
Set-TextValue "B12" "11.53"
Set-TextValue "C12" "39.58"

# Row 13 - Employment (absolute #)
Set-TextValue "B13" "730308.25"
Set-TextValue "C13" "2507698.75"

# Row 14 - Enterprises (% of total)
Set-TextValue "C14" "28.44"
Set-TextValue "D14" "98.14"
